# Actualización automática 2025-10-07 08:30:08
#
# Applies the "octubre" sales-update pass for GUERRERO FAREZ FABIAN MAURICIO
# across the three report sheets: VENTAS POR GRUPO, VENTA MENSUAL and
# CUMPLIMIENTO MENSUAL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" — per-client, per-product-group sales amounts.
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M11").Value = 1217.46
$wsGrupo.Range("M23").Value = 90.98
$wsGrupo.Range("M35").Value = 2587.7
$wsGrupo.Range("I42").Value = 86.40000000000001
$wsGrupo.Range("K42").Value = 406.08
$wsGrupo.Range("D46").Value = 457.92
$wsGrupo.Range("I48").Value = 81
$wsGrupo.Range("I53").Value = 26.1
$wsGrupo.Range("M53").Value = 42.77

# Row 56 is the "X de 54" fulfilment-count summary row; bump the counters for
# the columns that just received a first non-zero entry above.
$wsGrupo.Range("D56").Value = "1 de 54"
$wsGrupo.Range("I56").Value = "3 de 54"
$wsGrupo.Range("K56").Value = "1 de 54"
$wsGrupo.Range("M56").Value = "4 de 54"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" — per-client monthly sales, column F = octubre.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F11").Value = 1217.46
$wsMensual.Range("F23").Value = 90.98
$wsMensual.Range("F35").Value = 2587.7
$wsMensual.Range("F42").Value = 492.48
$wsMensual.Range("F46").Value = 457.92
$wsMensual.Range("F48").Value = 81
$wsMensual.Range("F55").Value = 41.22
$wsMensual.Range("F56").Value = 41.22
$wsMensual.Range("F60").Value = 4954.07

# ---------------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL" — budget vs. sales roll-up by product group.
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 457.92
$wsCumplimiento.Range("E3").Value = 17211.2270988183
$wsCumplimiento.Range("F3").Value = 0.0259163612957088

$wsCumplimiento.Range("D7").Value = 193.5
$wsCumplimiento.Range("E7").Value = 693.211016287574
$wsCumplimiento.Range("F7").Value = 0.2182221675897675

$wsCumplimiento.Range("D10").Value = 406.08
$wsCumplimiento.Range("E10").Value = 3474.99983534392
$wsCumplimiento.Range("F10").Value = 0.1046306742525474

$wsCumplimiento.Range("D12").Value = 3883
$wsCumplimiento.Range("E12").Value = 48780.12
$wsCumplimiento.Range("F12").Value = 0.07373281339958589

$wsCumplimiento.Range("D14").Value = 4912.85
$wsCumplimiento.Range("E14").Value = 94103.65661190613
$wsCumplimiento.Range("F14").Value = 0.04961647474855733

# The wider "VENTA"/"POR CUMPLIR" figures shrink columns E/F's auto-fit width
# by a notch (23->22, 28->25). Set the stored width directly (Excel's
# ColumnWidth property reports ~5/6 of a character narrower than what ends up
# persisted in the sheet's <col> width, hence the offset below).
$wsCumplimiento.Columns.Item(5).ColumnWidth = 21.166666666666668
$wsCumplimiento.Columns.Item(6).ColumnWidth = 24.166666666666668
